# "User Screen Code Updated"
#
# The Login sheet's test-data table loses its hyperlinked e-mail cells
# (and the "Hyperlink" cell-style that went with them), and a duplicate
# "User_TC004" row is collapsed into the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Remove the mailto: hyperlinks that lived on B2/B6/B7 (this also drops
# the <hyperlinks> block from the worksheet).
$ws.Hyperlinks.Delete()

# Strip the leftover "Hyperlink" look (underline, style index) from the
# cells that used to carry those hyperlinks.
$ws.Range("B2").Style = "Normal"
$ws.Range("B6").Style = "Normal"
$ws.Range("B7").Style = "Normal"

# Row 4 ("User_TC004") absorbs the outcome of the row below it, which is
# about to be deleted as a duplicate.
$ws.Cells.Item(4, 4).Value = "Login successful"

# Delete the now-redundant duplicate "User_TC004" row; rows 6 and 7 shift
# up to become rows 5 and 6.
$ws.Rows.Item(5).Delete()

# The "Hyperlink" named cell style is no longer used by any cell; drop it.
$wb.Styles.Item("Hyperlink").Delete()

# Reset the view back to the top-left of the (now smaller) sheet.
$ws.Range("A1").Select()
